# Switzerland Super League workbook update
# - Rows 97/98 (match ids 6811292 / 6811262) were swapped in the source feed.
# - Rows 101/102 (match ids 6810777 / 6811260) were swapped in the source feed.
# - Rows 177/178/179 received refreshed odds data (new ids/dates/teams/odds).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($Row, $Values) {
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

# ---- Row 97 (becomes what used to be row 98) ----
Set-RowValues 97 @{
    B  = 6811262
    F  = "Young Boys"
    G  = "St Gallen"
    H  = 3
    I  = 0
    J  = "H"
    K  = 2
    L  = 3.8
    M  = 3.3
    N  = 2.3
    O  = 4
    P  = 2.8
    Q  = -0.25
    R  = 2.025
    S  = 1.775
    T  = 3.5
    U  = 1.975
    V  = 1.875
    W  = 1.3
    X  = -1
    Y  = -1
    Z  = 1.025
    AA = -1
    AB = -1
    AC = 0.875
}

# ---- Row 98 (becomes what used to be row 97) ----
Set-RowValues 98 @{
    B  = 6811292
    F  = "Basel"
    G  = "Grasshoppers"
    H  = 0
    I  = 1
    J  = "A"
    K  = 2
    L  = 3.6
    M  = 3.5
    N  = 1.833
    O  = 3.8
    P  = 4.2
    Q  = -0.5
    R  = 1.825
    S  = 2.025
    T  = 2.75
    U  = 1.975
    V  = 1.875
    W  = -1
    X  = -1
    Y  = 3.2
    Z  = -1
    AA = 1.025
    AB = -1
    AC = 0.875
}

# ---- Row 101 (becomes what used to be row 102) ----
Set-RowValues 101 @{
    B  = 6811260
    F  = "FC Zurich"
    G  = "Lucerne"
    H  = 1
    I  = 1
    J  = "D"
    K  = 1.833
    L  = 3.6
    M  = 4.2
    N  = 1.7
    O  = 3.6
    P  = 4.75
    Q  = -0.75
    R  = 1.95
    S  = 1.9
    T  = 2.75
    U  = 1.975
    V  = 1.875
    W  = -1
    X  = 2.6
    Y  = -1
    Z  = -1
    AA = 0.8999999999999999
    AB = -1
    AC = 0.875
}

# ---- Row 102 (becomes what used to be row 101) ----
Set-RowValues 102 @{
    B  = 6810777
    F  = "Yverdon Sport FC"
    G  = "Stade LausanneOuchy"
    H  = 2
    I  = 1
    J  = "H"
    K  = 2.3
    L  = 3.6
    M  = 2.875
    N  = 3
    O  = 3.5
    P  = 2.375
    Q  = 0.25
    R  = 1.8
    S  = 2.05
    T  = 2.5
    U  = 1.925
    V  = 1.925
    W  = 2
    X  = -1
    Y  = -1
    Z  = 0.8
    AA = -1
    AB = 0.925
    AC = -1
}

# ---- Row 177 (refreshed odds, Stade LausanneOuchy vs Basel) ----
Set-RowValues 177 @{
    B = 7616909
    E = 45388.54166666666
    F = "Stade LausanneOuchy"
    G = "Basel"
    K = 3.5
    L = 3.6
    M = 2
    N = 3.5
    O = 3.6
    P = 2
    Q = 0.5
    R = 1.825
    S = 2.025
    T = 2.75
    U = 1.975
    V = 1.875
}

# ---- Row 178 (refreshed odds, Winterthur vs Lugano) ----
Set-RowValues 178 @{
    B = 7616935
    E = 45388.54166666666
    F = "Winterthur"
    G = "Lugano"
    K = 3.1
    L = 3.6
    M = 2.15
    N = 3.1
    O = 3.6
    P = 2.15
    Q = 0.25
    R = 1.925
    S = 1.925
    T = 2.75
    U = 1.9
    V = 1.95
}

# ---- Row 179 (refreshed odds, Servette vs FC Zurich) ----
Set-RowValues 179 @{
    B = 7616833
    E = 45388.64583333334
    F = "Servette"
    G = "FC Zurich"
    K = 1.85
    L = 3.6
    M = 3.8
    N = 1.85
    O = 3.6
    P = 3.8
    Q = -0.5
    R = 1.9
    S = 1.95
    T = 2.5
    U = 1.95
    V = 1.9
}
